$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.263.01"
$ws.Range("E2").Value = "  +1.99%  "

$ws.Range("D3").Value = "2.024.80"
$ws.Range("E3").Value = "  +3.72%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "247.32"
$ws.Range("E5").Value = "  +1.62%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.628"
$ws.Range("E6").Value = "  +0.74%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.58"
$ws.Range("E7").Value = "  +1.02%  "

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("E9").Value = "  +3.63%  "

$ws.Range("E10").Value = "  +3.18%  "

$ws.Range("E11").Value = "  +1.79%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.12"
$ws.Range("E12").Value = "  +6.95%  "

$ws.Range("D13").Value = "2.319.14"
$ws.Range("E13").Value = "  +3.52%  "

$ws.Range("E14").Value = "  +3.05%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.96"
$ws.Range("E15").Value = "  +1.90%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.46"
$ws.Range("E16").Value = "  +4.02%  "

$ws.Range("D17").Value = "2.025.49"
$ws.Range("E17").Value = "  +4.01%  "

$ws.Range("D18").Value = "37.193.81"
$ws.Range("E18").Value = "  +2.05%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "70.38"
$ws.Range("E19").Value = "  +1.71%  "

$ws.Range("D20").Value = "0.0₃0869"
$ws.Range("E20").Value = "  +2.21%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.24"
$ws.Range("E21").Value = "  +3.57%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "230.69"
$ws.Range("E22").Value = "  +0.68%  "

$ws.Range("E23").Value = "  -0.02%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.57"
$ws.Range("E24").Value = "  +4.93%  "

$ws.Range("E25").Value = "  -0.61%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.41"
$ws.Range("E26").Value = "  +2.99%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "163.55"
$ws.Range("E27").Value = "  +2.02%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.138"
$ws.Range("E28").Value = "  -3.69%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.79"
$ws.Range("E29").Value = "  +2.83%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.38"
$ws.Range("E30").Value = "  +6.07%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.122"
$ws.Range("E31").Value = "  +1.40%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0673"
$ws.Range("E32").Value = "  +10.42%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.78"
$ws.Range("E33").Value = "  +0.47%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.51"
$ws.Range("E34").Value = "  +11.00%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.48"
$ws.Range("E35").Value = "  +0.92%  "

$ws.Range("E36").Value = "  +5.42%  "

$ws.Range("E37").Value = "  +0.03%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.80"
$ws.Range("E38").Value = "  +1.70%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.37"
$ws.Range("E39").Value = "  -1.36%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.01"
$ws.Range("E40").Value = "  +3.23%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0976"
$ws.Range("E41").Value = "  +0.97%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0216"
$ws.Range("E42").Value = "  +2.92%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "16.96"
$ws.Range("E43").Value = "  +7.51%  "

$ws.Range("E44").Value = "  +1.76%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "91.56"
$ws.Range("E45").Value = "  +3.32%  "

$ws.Range("D46").Value = "1.380.49"
$ws.Range("E46").Value = "  +1.40%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.06"
$ws.Range("E47").Value = "  +3.32%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.45"
$ws.Range("E48").Value = "  +4.39%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.13"
$ws.Range("E49").Value = "  +15.75%  "

$ws.Range("E50").Value = "  +1.69%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "46.11"
$ws.Range("E51").Value = "  +1.95%  "
